$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" labels (column E) are put in reverse-chronological order
# (newest period first), and the matching "Valor Mora" (column F) values
# are updated to follow the new grouping: the 19 most-recent periods
# (rows 16-34) now carry the 31249 value while the remaining older
# periods (rows 35-50) carry 29509.

$periods = @(
    "2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
    "1712","1711","1710","1709","1708","1707","1706","1705"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    if ($row -le 34) {
        $ws.Range("F$row").Value = 31249
    } else {
        $ws.Range("F$row").Value = 29509
    }
}
